$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 186; everything currently at 186.. shifts down to 188..
$ws.Rows("186:187").Insert()

# Fill in the new row 186 (Primera quality) with the latest week's data
$ws.Range("A186").Value = 1
$ws.Range("B186").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C186").Value = "Arica y Parinacota"
$ws.Range("D186").Value = 44588
$ws.Range("E186").Value = 15
$ws.Range("F186").Value = 100112043
$ws.Range("G186").Value = "Pepino ensalada"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 120
$ws.Range("K186").Value = 5000
$ws.Range("L186").Value = 6000
$ws.Range("M186").Value = 5500
$ws.Range("N186").Value = "$/caja 70 unidades"
$ws.Range("O186").Value = "Región de Arica y Parinacota"
$ws.Range("P186").Value = 79
$ws.Range("Q186").Value = 70
$ws.Range("R186").Value = "Hortaliza"

# Fill in the new row 187 (Segunda quality) with the latest week's data
$ws.Range("A187").Value = 1
$ws.Range("B187").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C187").Value = "Arica y Parinacota"
$ws.Range("D187").Value = 44588
$ws.Range("E187").Value = 15
$ws.Range("F187").Value = 100112043
$ws.Range("G187").Value = "Pepino ensalada"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Segunda"
$ws.Range("J187").Value = 120
$ws.Range("K187").Value = 4000
$ws.Range("L187").Value = 5000
$ws.Range("M187").Value = 4500
$ws.Range("N187").Value = "$/caja 100 unidades"
$ws.Range("O187").Value = "Región de Arica y Parinacota"
$ws.Range("P187").Value = 45
$ws.Range("Q187").Value = 100
$ws.Range("R187").Value = "Hortaliza"
